$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 25,20
$data[0,0] = "ECs"
$data[0,1] = "Icam2"
$data[0,2] = "Itgal"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 37.78271433333333
$data[0,7] = 113.348143
$data[0,8] = 0.7852424281394679
$data[0,9] = 0.7883081272171703
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.1757273333333333
$data[0,13] = 0.527182
$data[0,14] = 0.00403052093865001
$data[0,15] = 0.004035374510811824
$data[0,16] = 6.639455635891778
$data[0,17] = 59.755100723026
$data[0,18] = 0.003164936048532501
$data[0,19] = 0.003181118523237973
$data[1,0] = "ECs"
$data[1,1] = "Icam2"
$data[1,2] = "Itgal"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 37.78271433333333
$data[1,7] = 113.348143
$data[1,8] = 0.7852424281394679
$data[1,9] = 0.7883081272171703
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.362774
$data[1,13] = 1.088322
$data[1,14] = 0.00832066460727691
$data[1,15] = 0.008330684390506021
$data[1,16] = 13.70658640956067
$data[1,17] = 123.359277686046
$data[1,18] = 0.006533738879952252
$data[1,19] = 0.006567146210317115
$data[2,0] = "ECs"
$data[2,1] = "Icam2"
$data[2,2] = "Itgal"
$data[2,3] = "M1"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 37.78271433333333
$data[2,7] = 113.348143
$data[2,8] = 0.7852424281394679
$data[2,9] = 0.7883081272171703
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 18.888837
$data[2,13] = 56.666511
$data[2,14] = 0.4332385383145499
$data[2,15] = 0.4337602461882952
$data[2,16] = 713.6715324598969
$data[2,17] = 6423.043792139073
$data[2,18] = 0.3401972817897111
$data[2,19] = 0.3419367273339537
$data[3,0] = "ECs"
$data[3,1] = "Icam2"
$data[3,2] = "Itgal"
$data[3,3] = "M2"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 37.78271433333333
$data[3,7] = 113.348143
$data[3,8] = 0.7852424281394679
$data[3,9] = 0.7883081272171703
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 24.01450566666666
$data[3,13] = 72.043517
$data[3,14] = 0.5508020072052686
$data[3,15] = 0.5514652855580013
$data[3,16] = 907.3332074598811
$data[3,17] = 8165.99886713893
$data[3,18] = 0.4325131055619578
$data[3,19] = 0.43472456648351
$data[4,0] = "ECs"
$data[4,1] = "Icam2"
$data[4,2] = "Itgal"
$data[4,3] = "sCs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 37.78271433333333
$data[4,7] = 113.348143
$data[4,8] = 0.7852424281394679
$data[4,9] = 0.7883081272171703
$data[4,10] = 2
$data[4,11] = 1
$data[4,12] = 0.1573175
$data[4,13] = 0.314635
$data[4,14] = 0.003608268934254619
$data[4,15] = 0.002408409352385472
$data[4,16] = 5.943882162134166
$data[4,17] = 35.663292972805
$data[4,18] = 0.002833365859314307
$data[4,19] = 0.001898568666151309
$data[5,0] = "FAPs"
$data[5,1] = "Icam2"
$data[5,2] = "Itgal"
$data[5,3] = "ECs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.270981
$data[5,7] = 3.812943
$data[5,8] = 0.02641494196933943
$data[5,9] = 0.02651806969185035
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.1757273333333333
$data[5,13] = 0.527182
$data[5,14] = 0.00403052093865001
$data[5,15] = 0.004035374510811824
$data[5,16] = 0.2233461018473333
$data[5,17] = 2.010114916626
$data[5,18] = 0.0001064659767006475
$data[5,19] = 0.0001070103425104244
$data[6,0] = "FAPs"
$data[6,1] = "Icam2"
$data[6,2] = "Itgal"
$data[6,3] = "FAPs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1.270981
$data[6,7] = 3.812943
$data[6,8] = 0.02641494196933943
$data[6,9] = 0.02651806969185035
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.362774
$data[6,13] = 1.088322
$data[6,14] = 0.00832066460727691
$data[6,15] = 0.008330684390506021
$data[6,16] = 0.4610788612939999
$data[6,17] = 4.149709751646
$data[6,18] = 0.000219789872747556
$data[6,19] = 0.0002209136692482485
$data[7,0] = "FAPs"
$data[7,1] = "Icam2"
$data[7,2] = "Itgal"
$data[7,3] = "M1"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.270981
$data[7,7] = 3.812943
$data[7,8] = 0.02641494196933943
$data[7,9] = 0.02651806969185035
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 18.888837
$data[7,13] = 56.666511
$data[7,14] = 0.4332385383145499
$data[7,15] = 0.4337602461882952
$data[7,16] = 24.007352939097
$data[7,17] = 216.066176451873
$data[7,18] = 0.01144397084846027
$data[7,19] = 0.01150248443797538
$data[8,0] = "FAPs"
$data[8,1] = "Icam2"
$data[8,2] = "Itgal"
$data[8,3] = "M2"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 1.270981
$data[8,7] = 3.812943
$data[8,8] = 0.02641494196933943
$data[8,9] = 0.02651806969185035
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 24.01450566666666
$data[8,13] = 72.043517
$data[8,14] = 0.5508020072052686
$data[8,15] = 0.5514652855580013
$data[8,16] = 30.52198042672566
$data[8,17] = 274.697823840531
$data[8,18] = 0.01454940305692285
$data[8,19] = 0.01462379487506323
$data[9,0] = "FAPs"
$data[9,1] = "Icam2"
$data[9,2] = "Itgal"
$data[9,3] = "sCs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 1.270981
$data[9,7] = 3.812943
$data[9,8] = 0.02641494196933943
$data[9,9] = 0.02651806969185035
$data[9,10] = 2
$data[9,11] = 1
$data[9,12] = 0.1573175
$data[9,13] = 0.314635
$data[9,14] = 0.003608268934254619
$data[9,15] = 0.002408409352385472
$data[9,16] = 0.1999475534675
$data[9,17] = 1.199685320805
$data[9,18] = 0.00009531221450810599
$data[9,19] = 0.00006386636705306211
$data[10,0] = "M1"
$data[10,1] = "Icam2"
$data[10,2] = "Itgal"
$data[10,3] = "ECs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 3.189763000000001
$data[10,7] = 9.569289000000001
$data[10,8] = 0.0662932054381191
$data[10,9] = 0.06655202362150626
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.1757273333333333
$data[10,13] = 0.527182
$data[10,14] = 0.00403052093865001
$data[10,15] = 0.004035374510811824
$data[10,16] = 0.5605285459553335
$data[10,17] = 5.044756913598001
$data[10,18] = 0.0002671961526085658
$data[10,19] = 0.0002685623397651728
$data[11,0] = "M1"
$data[11,1] = "Icam2"
$data[11,2] = "Itgal"
$data[11,3] = "FAPs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 3.189763000000001
$data[11,7] = 9.569289000000001
$data[11,8] = 0.0662932054381191
$data[11,9] = 0.06655202362150626
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.362774
$data[11,13] = 1.088322
$data[11,14] = 0.00832066460727691
$data[11,15] = 0.008330684390506021
$data[11,16] = 1.157163082562
$data[11,17] = 10.414467743058
$data[11,18] = 0.0005516035281918947
$data[11,19] = 0.0005544239043402703
$data[12,0] = "M1"
$data[12,1] = "Icam2"
$data[12,2] = "Itgal"
$data[12,3] = "M1"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 3.189763000000001
$data[12,7] = 9.569289000000001
$data[12,8] = 0.0662932054381191
$data[12,9] = 0.06655202362150626
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 18.888837
$data[12,13] = 56.666511
$data[12,14] = 0.4332385383145499
$data[12,15] = 0.4337602461882952
$data[12,16] = 60.25091337563101
$data[12,17] = 542.2582203806791
$data[12,18] = 0.0287207714241969
$data[12,19] = 0.02886762215039379
$data[13,0] = "M1"
$data[13,1] = "Icam2"
$data[13,2] = "Itgal"
$data[13,3] = "M2"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 3.189763000000001
$data[13,7] = 9.569289000000001
$data[13,8] = 0.0662932054381191
$data[13,9] = 0.06655202362150626
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 24.01450566666666
$data[13,13] = 72.043517
$data[13,14] = 0.5508020072052686
$data[13,15] = 0.5514652855580013
$data[13,16] = 76.60058163882367
$data[13,17] = 689.405234749413
$data[13,18] = 0.03651443061938723
$data[13,19] = 0.03670113071089681
$data[14,0] = "M1"
$data[14,1] = "Icam2"
$data[14,2] = "Itgal"
$data[14,3] = "sCs"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 3.189763000000001
$data[14,7] = 9.569289000000001
$data[14,8] = 0.0662932054381191
$data[14,9] = 0.06655202362150626
$data[14,10] = 2
$data[14,11] = 1
$data[14,12] = 0.1573175
$data[14,13] = 0.314635
$data[14,14] = 0.003608268934254619
$data[14,15] = 0.002408409352385472
$data[14,16] = 0.5018055407525001
$data[14,17] = 3.010833244515001
$data[14,18] = 0.0002392037137345245
$data[14,19] = 0.0001602845161102146
$data[15,0] = "M2"
$data[15,1] = "Icam2"
$data[15,2] = "Itgal"
$data[15,3] = "ECs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 5.311165
$data[15,7] = 15.933495
$data[15,8] = 0.1103825432989058
$data[15,9] = 0.1108134925816486
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 0.1757273333333333
$data[15,13] = 0.527182
$data[15,14] = 0.00403052093865001
$data[15,15] = 0.004035374510811824
$data[15,16] = 0.9333168623433334
$data[15,17] = 8.399851761090002
$data[15,18] = 0.0004448991520276814
$data[15,19] = 0.00044717394341802
$data[16,0] = "M2"
$data[16,1] = "Icam2"
$data[16,2] = "Itgal"
$data[16,3] = "FAPs"
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 5.311165
$data[16,7] = 15.933495
$data[16,8] = 0.1103825432989058
$data[16,9] = 0.1108134925816486
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 0.362774
$data[16,13] = 1.088322
$data[16,14] = 0.00832066460727691
$data[16,15] = 0.008330684390506021
$data[16,16] = 1.92675257171
$data[16,17] = 17.34077314539
$data[16,18] = 0.0009184561212884168
$data[16,19] = 0.0009231522329073951
$data[17,0] = "M2"
$data[17,1] = "Icam2"
$data[17,2] = "Itgal"
$data[17,3] = "M1"
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 5.311165
$data[17,7] = 15.933495
$data[17,8] = 0.1103825432989058
$data[17,9] = 0.1108134925816486
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 18.888837
$data[17,13] = 56.666511
$data[17,14] = 0.4332385383145499
$data[17,15] = 0.4337602461882952
$data[17,16] = 100.321729965105
$data[17,17] = 902.895569685945
$data[17,18] = 0.04782197171426049
$data[17,19] = 0.04806648782320073
$data[18,0] = "M2"
$data[18,1] = "Icam2"
$data[18,2] = "Itgal"
$data[18,3] = "M2"
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = 5.311165
$data[18,7] = 15.933495
$data[18,8] = 0.1103825432989058
$data[18,9] = 0.1108134925816486
$data[18,10] = 3
$data[18,11] = 1
$data[18,12] = 24.01450566666666
$data[18,13] = 72.043517
$data[18,14] = 0.5508020072052686
$data[18,15] = 0.5514652855580013
$data[18,16] = 127.5450019891016
$data[18,17] = 1147.905017901915
$data[18,18] = 0.06079892640945981
$data[18,19] = 0.06110979433021833
$data[19,0] = "M2"
$data[19,1] = "Icam2"
$data[19,2] = "Itgal"
$data[19,3] = "sCs"
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 5.311165
$data[19,7] = 15.933495
$data[19,8] = 0.1103825432989058
$data[19,9] = 0.1108134925816486
$data[19,10] = 2
$data[19,11] = 1
$data[19,12] = 0.1573175
$data[19,13] = 0.314635
$data[19,14] = 0.003608268934254619
$data[19,15] = 0.002408409352385472
$data[19,16] = 0.8355391998875
$data[19,17] = 5.013235199325
$data[19,18] = 0.0003982899018694573
$data[19,19] = 0.0002668842519041407
$data[20,0] = "sCs"
$data[20,1] = "Icam2"
$data[20,2] = "Itgal"
$data[20,3] = "ECs"
$data[20,4] = 1
$data[20,5] = 0.5
$data[20,6] = 0.5613635
$data[20,7] = 1.122727
$data[20,8] = 0.01166688115416775
$data[20,9] = 0.007808286887824462
$data[20,10] = 3
$data[20,11] = 1
$data[20,12] = 0.1757273333333333
$data[20,13] = 0.527182
$data[20,14] = 0.00403052093865001
$data[20,15] = 0.004035374510811824
$data[20,16] = 0.09864691088566668
$data[20,17] = 0.5918814653140001
$data[20,18] = 0.00004702360878061429
$data[20,19] = 0.00003150936188023301
$data[21,0] = "sCs"
$data[21,1] = "Icam2"
$data[21,2] = "Itgal"
$data[21,3] = "FAPs"
$data[21,4] = 1
$data[21,5] = 0.5
$data[21,6] = 0.5613635
$data[21,7] = 1.122727
$data[21,8] = 0.01166688115416775
$data[21,9] = 0.007808286887824462
$data[21,10] = 3
$data[21,11] = 1
$data[21,12] = 0.362774
$data[21,13] = 1.088322
$data[21,14] = 0.00832066460727691
$data[21,15] = 0.008330684390506021
$data[21,16] = 0.203648082349
$data[21,17] = 1.221888494094
$data[21,18] = 0.00009707620509678953
$data[21,19] = 0.00006504837369299208
$data[22,0] = "sCs"
$data[22,1] = "Icam2"
$data[22,2] = "Itgal"
$data[22,3] = "M1"
$data[22,4] = 1
$data[22,5] = 0.5
$data[22,6] = 0.5613635
$data[22,7] = 1.122727
$data[22,8] = 0.01166688115416775
$data[22,9] = 0.007808286887824462
$data[22,10] = 3
$data[22,11] = 1
$data[22,12] = 18.888837
$data[22,13] = 56.666511
$data[22,14] = 0.4332385383145499
$data[22,15] = 0.4337602461882952
$data[22,16] = 10.6035036492495
$data[22,17] = 63.621021895497
$data[22,18] = 0.005054542537921204
$data[22,19] = 0.003386924442771576
$data[23,0] = "sCs"
$data[23,1] = "Icam2"
$data[23,2] = "Itgal"
$data[23,3] = "M2"
$data[23,4] = 1
$data[23,5] = 0.5
$data[23,6] = 0.5613635
$data[23,7] = 1.122727
$data[23,8] = 0.01166688115416775
$data[23,9] = 0.007808286887824462
$data[23,10] = 3
$data[23,11] = 1
$data[23,12] = 24.01450566666666
$data[23,13] = 72.043517
$data[23,14] = 0.5508020072052686
$data[23,15] = 0.5514652855580013
$data[23,16] = 13.48086695180983
$data[23,17] = 80.885201710859
$data[23,18] = 0.006426141557540915
$data[23,19] = 0.004305999158312914
$data[24,0] = "sCs"
$data[24,1] = "Icam2"
$data[24,2] = "Itgal"
$data[24,3] = "sCs"
$data[24,4] = 1
$data[24,5] = 0.5
$data[24,6] = 0.5613635
$data[24,7] = 1.122727
$data[24,8] = 0.01166688115416775
$data[24,9] = 0.007808286887824462
$data[24,10] = 2
$data[24,11] = 1
$data[24,12] = 0.1573175
$data[24,13] = 0.314635
$data[24,14] = 0.003608268934254619
$data[24,15] = 0.002408409352385472
$data[24,16] = 0.08831230241125
$data[24,17] = 0.353249209645
$data[24,18] = 0.00004209724482822415
$data[24,19] = 0.00001880555116674529

$ws.Range("A2:T26").Value = $data

Write-Host "done"
